$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-34 with the new date/value pairs
$data = @(
    @("05/05/2021 01:39:46", 13.99),
    @("05/05/2021 01:40:48", 12.964),
    @("05/05/2021 01:41:50", 13.814),
    @("05/05/2021 01:42:55", 13.814),
    @("05/05/2021 01:44:00", 13.862),
    @("05/05/2021 01:45:04", 54.964),
    @("05/05/2021 01:46:09", 66.68300000000001),
    @("05/05/2021 01:51:31", 133.993),
    @("05/05/2021 01:52:37", 79.185),
    @("05/05/2021 01:53:43", 189.745),
    @("05/05/2021 01:54:50", 66.16),
    @("05/05/2021 02:00:08", 65.178),
    @("05/05/2021 02:01:17", 111.654),
    @("05/05/2021 02:06:33", 622.159),
    @("05/05/2021 02:07:38", 25.656),
    @("05/05/2021 02:12:49", 12.944),
    @("05/05/2021 02:13:53", 13.124),
    @("05/05/2021 02:14:57", 15.232),
    @("05/05/2021 02:16:00", 12.916),
    @("05/05/2021 02:17:04", 12.864),
    @("05/05/2021 02:18:07", 15.361),
    @("05/05/2021 02:19:12", 13.726),
    @("05/05/2021 02:20:16", 9.776),
    @("05/05/2021 02:21:21", 63.891),
    @("05/05/2021 02:21:31", 63.891),
    @("05/05/2021 02:21:40", 63.891),
    @("05/05/2021 02:21:48", 63.891),
    @("05/05/2021 02:21:58", 113.958),
    @("05/05/2021 02:22:08", 113.958),
    @("05/05/2021 02:22:18", 113.958),
    @("05/05/2021 02:22:30", 113.958),
    @("05/05/2021 02:22:40", 113.958),
    @("05/05/2021 02:22:49", 115.033)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Remove the now-unused trailing rows (old rows 35-53)
$ws.Range("A35:B53").Clear()
